$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.604577898979187
$ws.Range("B1").Value = 1.486402750015259
$ws.Range("C1").Value = 4.741234302520752
$ws.Range("D1").Value = 1.403598427772522
$ws.Range("E1").Value = 0.6569319367408752
